$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("U3").Value = 4.5
$ws.Range("V3").Value = 1.21
$ws.Range("S4").Value = 3.1
$ws.Range("T4").Value = 1.36
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 4.75
$ws.Range("I6").Value = 7
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("Q6").Value = 1.45
$ws.Range("R6").Value = 2.75
$ws.Range("S6").Value = 1.93
$ws.Range("T6").Value = 1.93
$ws.Range("U6").Value = 2.6
$ws.Range("V6").Value = 1.51
$ws.Range("AI6").Value = 11
$ws.Range("AJ6").Value = 9.5
$ws.Range("AM6").Value = 15
$ws.Range("AN6").Value = 34
$ws.Range("AO6").Value = 21
$ws.Range("AP6").Value = 81
$ws.Range("H7").Value = 4.2
$ws.Range("K7").Value = 2.38
$ws.Range("U7").Value = 2.03
$ws.Range("V7").Value = 1.83
$ws.Range("W7").Value = 2.5
$ws.Range("X7").Value = 1.53
$ws.Range("Y7").Value = 1.3
$ws.Range("Z7").Value = 3.4
$ws.Range("AA7").Value = 1.62
$ws.Range("AB7").Value = 2.2
$ws.Range("AD7").Value = 10
$ws.Range("AM7").Value = 15
$ws.Range("AR7").Value = 29
$ws.Range("G8").Value = 2.1
$ws.Range("I8").Value = 3.2
$ws.Range("J8").Value = 2.75
$ws.Range("G9").Value = 1.55
$ws.Range("H9").Value = 4.1
$ws.Range("J9").Value = 2.1
$ws.Range("K9").Value = 2.4
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 4.33
$ws.Range("S9").Value = 1.67
$ws.Range("T9").Value = 2.15
$ws.Range("U9").Value = 2.03
$ws.Range("V9").Value = 1.78
$ws.Range("W9").Value = 2.63
$ws.Range("X9").Value = 1.44
$ws.Range("Y9").Value = 1.3
$ws.Range("Z9").Value = 3.4
$ws.Range("AA9").Value = 1.73
$ws.Range("AB9").Value = 2
$ws.Range("AC9").Value = 8
$ws.Range("AD9").Value = 8
$ws.Range("AF9").Value = 11
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 23
$ws.Range("AI9").Value = 13
$ws.Range("AJ9").Value = 8
$ws.Range("AM9").Value = 17
$ws.Range("AN9").Value = 29
$ws.Range("AO9").Value = 17
$ws.Range("AS9").Value = 201
$ws.Range("G14").Value = 3.4
$ws.Range("I14").Value = 2.05
$ws.Range("J14").Value = 4
$ws.Range("M14").Value = 1.04
$ws.Range("N14").Value = 13
$ws.Range("O14").Value = 1.25
$ws.Range("P14").Value = 3.75
$ws.Range("S14").Value = 1.88
$ws.Range("T14").Value = 1.98
$ws.Range("Y14").Value = 1.4
$ws.Range("Z14").Value = 2.75
$ws.Range("AA14").Value = 1.75
$ws.Range("AB14").Value = 2
$ws.Range("AF14").Value = 41
$ws.Range("AH14").Value = 34
$ws.Range("AI14").Value = 11
$ws.Range("AK14").Value = 15
$ws.Range("AL14").Value = 51
$ws.Range("AM14").Value = 8
$ws.Range("AN14").Value = 10
$ws.Range("AR14").Value = 26
$ws.Range("S15").Value = 2.25
$ws.Range("T15").Value = 1.62
$ws.Range("W15").Value = 4
$ws.Range("X15").Value = 1.22
$ws.Range("AS15").Value = 1250
$ws.Range("G19").Value = 2.7
$ws.Range("I19").Value = 2.9
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 7
$ws.Range("AO19").Value = 12
$ws.Range("G20").Value = 5.6
$ws.Range("H20").Value = 4.15
$ws.Range("J20").Value = 5.5
$ws.Range("K20").Value = 2.32
$ws.Range("L20").Value = 2
$ws.Range("O20").Value = 1.21
$ws.Range("T20").Value = 2.12
$ws.Range("W20").Value = 2.52
$ws.Range("X20").Value = 1.47
$ws.Range("AA20").Value = 1.75
$ws.Range("AB20").Value = 1.98
$ws.Range("AD20").Value = 35
$ws.Range("AE20").Value = 17.5
$ws.Range("AF20").Value = 110
$ws.Range("AG20").Value = 55
$ws.Range("AK20").Value = 16
$ws.Range("AL20").Value = 65
$ws.Range("AM20").Value = 7.9
$ws.Range("AN20").Value = 7.8
$ws.Range("AO20").Value = 8
$ws.Range("AP20").Value = 11
$ws.Range("AQ20").Value = 11.25
$ws.Range("AR20").Value = 22
$ws.Range("AS20").Value = 450
